# Generate Report for Handoff
# A new handoff batch was produced: every file currently awaiting handoff
# (status "Ready for handoff") as well as the file whose handback transform
# failed gets its "Latest Handoff Datetime" column (D) refreshed to the
# timestamp of this new handoff run - one timestamp per locale sheet.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

$ws_zh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $ws_zh.Cells.Item($r, 4).Value = "2016-03-08 22:33:03"
}

$ws_de = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $ws_de.Cells.Item($r, 4).Value = "2016-03-08 22:33:13"
}
